$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.716.73"
$ws.Range("E2").Value = "  -4.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.017.83"
$ws.Range("E3").Value = "  -6.25%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.46"
$ws.Range("E5").Value = "  -2.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.22"
$ws.Range("E6").Value = "  -7.64%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.012.72"
$ws.Range("E8").Value = "  -6.32%  "
$ws.Range("E9").Value = "  -3.04%  "
$ws.Range("E10").Value = "  -7.20%  "
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.441"
$ws.Range("E13").Value = "  -7.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.78"
$ws.Range("E14").Value = "  -8.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.119"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.513.58"
$ws.Range("E16").Value = "  -6.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.014.95"
$ws.Range("E17").Value = "  -6.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.697.90"
$ws.Range("E18").Value = "  -5.05%  "
$ws.Range("E19").Value = "  -2.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "434.54"
$ws.Range("E20").Value = "  -6.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.14"
$ws.Range("E21").Value = "  -6.67%  "
$ws.Range("E22").Value = "  -5.56%  "
$ws.Range("E23").Value = "  -8.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.86"
$ws.Range("E24").Value = "  -4.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.60"
$ws.Range("E25").Value = "  -4.96%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.57"
$ws.Range("E28").Value = "  -4.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.35"
$ws.Range("E29").Value = "  -7.12%  "
$ws.Range("E30").Value = "  -7.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.20"
$ws.Range("E31").Value = "  -10.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.45"
$ws.Range("E32").Value = "  -7.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0937"
$ws.Range("E33").Value = "  -9.38%  "
$ws.Range("E34").Value = "  -12.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.959"
$ws.Range("E35").Value = "  -8.11%  "
$ws.Range("E36").Value = "  -5.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "50.59"
$ws.Range("E37").Value = "  -2.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0668"
$ws.Range("E38").Value = "  -9.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.52"
$ws.Range("E39").Value = "  +4.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0360"
$ws.Range("E40").Value = "  -8.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "388.88"
$ws.Range("E41").Value = "  -4.15%  "
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("E43").Value = "  -9.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.665.78"
$ws.Range("E44").Value = "  -6.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.236"
$ws.Range("E46").Value = "  -7.99%  "
$ws.Range("E47").Value = "  -6.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.09"
$ws.Range("E48").Value = "  -7.73%  "
$ws.Range("E49").Value = "  -3.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.73"
$ws.Range("E50").Value = "  -8.41%  "
$ws.Range("E51").Value = "  +2.09%  "
